$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.006.50"
$ws.Range("E2").Value = "  +2.04%  "

$ws.Range("D3").Value = "3.283.76"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.28%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  +1.77%  "

$ws.Range("D9").Value = "3.277.80"
$ws.Range("E9").Value = "  -0.74%  "

$ws.Range("E10").Value = "  -2.11%  "

$ws.Range("E11").Value = "  -0.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000269"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "694.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.22%  "

$ws.Range("D15").Value = "3.811.21"
$ws.Range("E15").Value = "  -0.61%  "

$ws.Range("E16").Value = "  -1.82%  "

$ws.Range("D17").Value = "67.089.78"
$ws.Range("E17").Value = "  +1.91%  "

$ws.Range("E18").Value = "  +0.98%  "

$ws.Range("D19").Value = "3.284.10"
$ws.Range("E19").Value = "  -0.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.888"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.30%  "

$ws.Range("E23").Value = "  -4.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.55%  "

$ws.Range("E26").Value = "  -2.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "570.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.95%  "

$ws.Range("D33").Value = "3.890.87"
$ws.Range("E33").Value = "  +1.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.96%  "

$ws.Range("E35").Value = "  -1.84%  "

$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("E38").Value = "  -11.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.128"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.70%  "

$ws.Range("E40").Value = "  +0.37%  "

$ws.Range("E41").Value = "  -1.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "31.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.38%  "

$ws.Range("D43").Value = "0.0₃0672"
$ws.Range("E43").Value = "  -4.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.47%  "

$ws.Range("E45").Value = "  -1.52%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0406"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.127"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.36%  "

$ws.Range("E48").Value = "  +0.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.35%  "

$ws.Range("E50").Value = "  +7.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "130.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.21%  "
